# Updated cryptos list - applies per-cell text updates while preserving cell
# type/style (force Text number format so Excel doesn't coerce numeric-looking
# strings like "0.999" / "4.30" into actual numbers, then restore the
# original style so no new cell formatting is introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@('D2', '58.638.88')
    ,@('E2', '  -4.81%  ')
    ,@('D3', '2.468.32')
    ,@('E3', '  -4.35%  ')
    ,@('D4', '0.999')
    ,@('E4', '  -0.07%  ')
    ,@('D5', '532.62')
    ,@('E5', '  -3.80%  ')
    ,@('D6', '144.07')
    ,@('E6', '  -6.88%  ')
    ,@('D7', '0.997')
    ,@('D8', '0.567')
    ,@('E8', '  -4.87%  ')
    ,@('D9', '2.491.34')
    ,@('E9', '  -3.51%  ')
    ,@('D10', '0.0993')
    ,@('E10', '  -4.69%  ')
    ,@('E11', '  -2.16%  ')
    ,@('D12', '5.56')
    ,@('E12', '  +1.75%  ')
    ,@('D13', '0.351')
    ,@('E13', '  -3.49%  ')
    ,@('D14', '2.899.98')
    ,@('E14', '  -4.45%  ')
    ,@('D15', '23.72')
    ,@('E15', '  -6.69%  ')
    ,@('D16', '58.525.95')
    ,@('E16', '  -4.87%  ')
    ,@('E17', '  -4.55%  ')
    ,@('D18', '2.471.23')
    ,@('E18', '  -4.30%  ')
    ,@('D19', '11.27')
    ,@('E19', '  -2.65%  ')
    ,@('D20', '4.30')
    ,@('E20', '  -5.21%  ')
    ,@('D21', '322.26')
    ,@('E21', '  -4.70%  ')
    ,@('E22', '  -0.27%  ')
    ,@('D23', '5.72')
    ,@('E23', '  -5.31%  ')
    ,@('D24', '60.62')
    ,@('E24', '  -3.43%  ')
    ,@('D25', '0.437')
    ,@('E25', '  -11.49%  ')
    ,@('D26', '0.995')
    ,@('E26', '  -0.21%  ')
    ,@('E27', '  -4.67%  ')
    ,@('D28', '2.580.68')
    ,@('E28', '  -4.53%  ')
    ,@('D29', '7.68')
    ,@('E29', '  -4.63%  ')
    ,@('D30', '6.93')
    ,@('E30', '  -1.77%  ')
    ,@('D31', '0.0₃0771')
    ,@('E31', '  -8.05%  ')
    ,@('D32', '1.78')
    ,@('E32', '  -7.28%  ')
    ,@('E33', '  -6.02%  ')
    ,@('D34', '0.996')
    ,@('E34', '  -0.30%  ')
    ,@('D35', '157.90')
    ,@('E35', '  -1.04%  ')
    ,@('D36', '1.40')
    ,@('E36', '  -0.95%  ')
    ,@('D37', '18.47')
    ,@('E37', '  -3.76%  ')
    ,@('D38', '4.39')
    ,@('E38', '  -6.03%  ')
    ,@('E39', '  -10.87%  ')
    ,@('D40', '5.76')
    ,@('E40', '  -4.69%  ')
    ,@('D41', '304.75')
    ,@('E41', '  -9.38%  ')
    ,@('D42', '36.51')
    ,@('E42', '  -2.37%  ')
    ,@('D43', '3.70')
    ,@('E43', '  -5.95%  ')
    ,@('D44', '0.807')
    ,@('E44', '  -9.53%  ')
    ,@('D45', '0.996')
    ,@('E45', '  -0.16%  ')
    ,@('D47', '0.591')
    ,@('E47', '  -2.58%  ')
    ,@('D48', '124.16')
    ,@('E48', '  -0.05%  ')
    ,@('D49', '0.0922')
    ,@('E49', '  -4.60%  ')
    ,@('D50', '0.0519')
    ,@('E50', '  -4.83%  ')
    ,@('B51', 'EnergySwap')
    ,@('C51', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens')
    ,@('D51', '18.47')
    ,@('E51', '  -5.54%  ')
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

